# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the refreshed scrape output (commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of row -> new F value for the "展览" sheet (sheet index 1)
$sheetExhibitionUpdates = @{
    3  = 160
    5  = 10
    6  = 529
    7  = 1589
    10 = 1326
    12 = 11
    13 = 231
    14 = 173
    16 = 10
    18 = 241
    19 = 133
    20 = 201
    21 = 188
}

# Map of row -> new F value for the "全部类型" sheet (sheet index 4)
$sheetAllTypesUpdates = @{
    3  = 160
    5  = 10
    6  = 529
    7  = 1589
    11 = 1326
    13 = 11
    14 = 231
    15 = 173
    17 = 10
    19 = 241
    20 = 133
    21 = 201
    22 = 188
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $sheetExhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $sheetExhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheetAllTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $sheetAllTypesUpdates[$row]
}
